$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = [double]"3.0809169402346015E-4"
$ws.Range("A3").Value = [double]"1.701587316347286E-4"
$ws.Range("H3").Value = [double]"4.768362045288086"
$ws.Range("A4").Value = [double]"1.379329478368163E-4"
$ws.Range("H4").Value = [double]"4.698826789855957"
